# Add a new "Greece" market test-data sheet, cloned from the existing
# "Croatia" sheet, positioned right after it.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Copy Croatia so the new sheet inherits all formatting/layout, placing
# the copy immediately after Croatia.
$croatia.Copy($null, $croatia)

# The newly created sheet becomes the active sheet, named "Croatia (2)".
$greece = $wb.ActiveSheet
$greece.Name = "Greece"

# Update the market-specific cells.
$greece.Range("B1").Value = "Greece Market"
$greece.Range("B3").Value = "NGC-4119/T3189"

# Croatia tab should no longer be the selected/active tab, and its
# selection reverts to a full-column selection (A1:XFD1048576).
$croatia.Select()
$croatia.Range("A1:XFD1048576").Select()

# Greece becomes the active/selected tab with B3 selected.
$greece.Select()
$greece.Range("B3").Select()
